$wb = $excel.ActiveWorkbook

# --- Sheet "Matriz": update matrix values ---
$ws1 = $wb.Worksheets.Item("Matriz")

$ws1.Range("B2").Value = 2
$ws1.Range("C2").Value = 2

$ws1.Range("A3").Value = 3
$ws1.Range("B3").Value = 6
$ws1.Range("C3").Value = 3

$ws1.Range("A4").Value = 6
$ws1.Range("B4").Value = 6
$ws1.Range("C4").Value = 9

# Wrap text on the cells that changed style (B3, A4, B4)
$ws1.Range("B3").WrapText = $true
$ws1.Range("A4").WrapText = $true
$ws1.Range("B4").WrapText = $true

# Move the (non-active) selection on this sheet to A4
[void]$ws1.Range("A4").Select()

# --- Sheet "Interacao": add header row with labels ---
$ws2 = $wb.Worksheets.Item("Interacao")
$ws2.Range("A1").Value = "cx"
$ws2.Range("B1").Value = "cy"
$ws2.Range("C1").Value = "cz"

# Restore "Interacao" as the active sheet (it was active before the edit)
$ws2.Activate()
